$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-10 07:18:26"
$ws.Range("G2").Value = "113 cm"
$ws.Range("E3").Value = "2026-02-10 07:18:28"
$ws.Range("G3").Value = "186 cm"
$ws.Range("I3").Value = "8.6 mm"
$ws.Range("E4").Value = "2026-02-10 07:18:31"
$ws.Range("J4").Value = "1005.0 hPa"
$ws.Range("E5").Value = "2026-02-10 07:18:33"
$ws.Range("I5").Value = "13.6 mm"
$ws.Range("E6").Value = "2026-02-10 07:18:35"
$ws.Range("J6").Value = "1005.1 hPa"
$ws.Range("N6").Value = "6.0 °C 6:36 TU"
$ws.Range("O6").Value = "7.2 °C"
$ws.Range("E7").Value = "2026-02-10 07:18:37"
$ws.Range("H7").Value = "'78%"
$ws.Range("J7").Value = "1005.2 hPa"
$ws.Range("E8").Value = "2026-02-10 07:18:40"
$ws.Range("J8").Value = "1005.2 hPa"
$ws.Range("O8").Value = "8.6 °C"
$ws.Range("E9").Value = "2026-02-10 07:18:42"
$ws.Range("I9").Value = "0.2 mm"
$ws.Range("N9").Value = "4.9 °C 6:53 TU"
$ws.Range("O9").Value = "6.3 °C"
$ws.Range("E10").Value = "2026-02-10 07:18:45"
$ws.Range("N10").Value = "4.5 °C 6:44 TU"
$ws.Range("O10").Value = "6.7 °C"
$ws.Range("E11").Value = "2026-02-10 07:18:47"
$ws.Range("E12").Value = "2026-02-10 07:18:49"
$ws.Range("N12").Value = "5.0 °C 6:36 TU"
$ws.Range("O12").Value = "6.6 °C"
$ws.Range("E13").Value = "2026-02-10 07:18:51"
$ws.Range("I13").Value = "2.5 mm"
$ws.Range("J13").Value = "1008.4 hPa"
$ws.Range("N13").Value = "2.3 °C 6:58 TU"
$ws.Range("O13").Value = "2.5 °C"
$ws.Range("E14").Value = "2026-02-10 07:18:54"
$ws.Range("E15").Value = "2026-02-10 07:18:56"
$ws.Range("E16").Value = "2026-02-10 07:18:59"
$ws.Range("G16").Value = "80 cm"
$ws.Range("I16").Value = "13.3 mm"
$ws.Range("M16").Value = "0.5 °C 6:59 TU"
$ws.Range("O16").Value = "-0.3 °C"
$ws.Range("E17").Value = "2026-02-10 07:19:01"
$ws.Range("H17").Value = "'89%"
$ws.Range("O17").Value = "3.0 °C"
$ws.Range("E18").Value = "2026-02-10 07:19:04"
$ws.Range("N18").Value = "4.5 °C 6:38 TU"
$ws.Range("O18").Value = "6.9 °C"
$ws.Range("E19").Value = "2026-02-10 07:19:06"
$ws.Range("N19").Value = "3.1 °C 6:34 TU"
$ws.Range("O19").Value = "3.8 °C"
$ws.Range("E20").Value = "2026-02-10 07:19:08"
$ws.Range("M20").Value = "-0.3 °C 6:59 TU"
$ws.Range("O20").Value = "-1.0 °C"
$ws.Range("E21").Value = "2026-02-10 07:19:11"
$ws.Range("I21").Value = "3.7 mm"
$ws.Range("J21").Value = "1007.7 hPa"
$ws.Range("E22").Value = "2026-02-10 07:19:13"
$ws.Range("G22").Value = "126 cm"
$ws.Range("I22").Value = "0.4 mm"
$ws.Range("M22").Value = "-1.2 °C 6:59 TU"
$ws.Range("E23").Value = "2026-02-10 07:19:15"
$ws.Range("G23").Value = "185 cm"
$ws.Range("I23").Value = "11.6 mm"
$ws.Range("E24").Value = "2026-02-10 07:19:18"
$ws.Range("J24").Value = "1007.1 hPa"
$ws.Range("N24").Value = "8.0 °C 6:30 TU"
$ws.Range("E25").Value = "2026-02-10 07:19:20"
$ws.Range("G25").Value = "120 cm"
$ws.Range("H25").Value = "'95%"
$ws.Range("I25").Value = "7.5 mm"
$ws.Range("E26").Value = "2026-02-10 07:19:23"
$ws.Range("O26").Value = "3.4 °C"
$ws.Range("E27").Value = "2026-02-10 07:19:25"
$ws.Range("H27").Value = "'98%"
$ws.Range("I27").Value = "2.1 mm"
$ws.Range("L27").Value = "34.6 km/h - 254º 6:59 TU"
$ws.Range("M27").Value = "0.2 °C 6:59 TU"
$ws.Range("O27").Value = "-0.3 °C"
$ws.Range("E28").Value = "2026-02-10 07:19:28"
$ws.Range("J28").Value = "1005.7 hPa"
$ws.Range("N28").Value = "3.2 °C 6:55 TU"
$ws.Range("O28").Value = "4.8 °C"
$ws.Range("E29").Value = "2026-02-10 07:19:30"
$ws.Range("H29").Value = "'94%"
$ws.Range("O29").Value = "8.3 °C"
$ws.Range("E30").Value = "2026-02-10 07:19:33"
$ws.Range("J30").Value = "1005.1 hPa"
$ws.Range("E31").Value = "2026-02-10 07:19:35"
$ws.Range("E32").Value = "2026-02-10 07:19:37"
$ws.Range("M32").Value = "8.4 °C 6:56 TU"
$ws.Range("E33").Value = "2026-02-10 07:19:40"
$ws.Range("I33").Value = "6.0 mm"
$ws.Range("E34").Value = "2026-02-10 07:19:42"
$ws.Range("I34").Value = "3.2 mm"
$ws.Range("N34").Value = "0.9 °C 6:39 TU"
$ws.Range("E35").Value = "2026-02-10 07:19:45"
$ws.Range("H35").Value = "'82%"
$ws.Range("J35").Value = "1005.7 hPa"
$ws.Range("M35").Value = "11.2 °C 6:59 TU"
$ws.Range("E36").Value = "2026-02-10 07:19:47"
$ws.Range("H36").Value = "'95%"
$ws.Range("N36").Value = "6.2 °C 6:33 TU"
$ws.Range("O36").Value = "8.5 °C"
$ws.Range("E37").Value = "2026-02-10 07:19:50"
$ws.Range("J37").Value = "1007.2 hPa"
$ws.Range("E38").Value = "2026-02-10 07:19:52"
$ws.Range("N38").Value = "6.0 °C 6:59 TU"
$ws.Range("O38").Value = "7.5 °C"
$ws.Range("E39").Value = "2026-02-10 07:19:54"
$ws.Range("I39").Value = "3.1 mm"
$ws.Range("M39").Value = "0.9 °C 6:52 TU"
$ws.Range("O39").Value = "-0.1 °C"
$ws.Range("E40").Value = "2026-02-10 07:19:57"
$ws.Range("I40").Value = "4.1 mm"
$ws.Range("J40").Value = "1008.5 hPa"
$ws.Range("E41").Value = "2026-02-10 07:19:59"
$ws.Range("E42").Value = "2026-02-10 07:20:02"
$ws.Range("H42").Value = "'100%"
$ws.Range("N42").Value = "6.4 °C 6:59 TU"
$ws.Range("O42").Value = "7.8 °C"
$ws.Range("E43").Value = "2026-02-10 07:20:04"
$ws.Range("N43").Value = "5.0 °C 6:58 TU"
$ws.Range("E44").Value = "2026-02-10 07:20:06"
$ws.Range("G44").Value = "219 cm"
$ws.Range("I44").Value = "7.8 mm"
$ws.Range("E45").Value = "2026-02-10 07:20:08"
$ws.Range("I45").Value = "18.0 mm"
$ws.Range("J45").Value = "1008.1 hPa"
$ws.Range("E46").Value = "2026-02-10 07:20:11"
$ws.Range("J46").Value = "1006.9 hPa"
